$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Retrocue column (G) values for the practice session rows.
# Cyclic relabeling: REV_BTH -> REP_1ST -> REP_2ND -> REV_BTH
$ws.Range("G5").Value  = "REP_1ST"
$ws.Range("G6").Value  = "REP_1ST"
$ws.Range("G7").Value  = "REP_1ST"

$ws.Range("G8").Value  = "REP_2ND"
$ws.Range("G9").Value  = "REP_2ND"
$ws.Range("G10").Value = "REP_2ND"

$ws.Range("G11").Value = "REV_BTH"
$ws.Range("G12").Value = "REV_BTH"
$ws.Range("G13").Value = "REV_BTH"

$ws.Range("G18").Value = "REP_1ST"
$ws.Range("G19").Value = "REP_2ND"
$ws.Range("G20").Value = "REV_BTH"

# Move the active cell selection to L17, matching the final saved state.
$ws.Range("L17").Select()
